$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.814.09"
$ws.Range("E2").Value = "  -5.68%  "

$ws.Range("D3").Value = "2.970.69"
$ws.Range("E3").Value = "  -6.36%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.54"
$ws.Range("E5").Value = "  -6.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "123.24"
$ws.Range("E6").Value = "  -9.29%  "

$ws.Range("E7").Value = "  +0.22%  "

$ws.Range("D8").Value = "2.960.60"
$ws.Range("E8").Value = "  -6.62%  "

$ws.Range("E9").Value = "  -2.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.129"
$ws.Range("E10").Value = "  -9.75%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.87"
$ws.Range("E11").Value = "  -9.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.438"
$ws.Range("E12").Value = "  -3.81%  "

$ws.Range("E13").Value = "  -9.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.07"
$ws.Range("E14").Value = "  -8.00%  "

$ws.Range("E15").Value = "  +0.05%  "

$ws.Range("D16").Value = "3.461.90"
$ws.Range("E16").Value = "  -6.30%  "

$ws.Range("D17").Value = "2.970.07"
$ws.Range("E17").Value = "  -6.37%  "

$ws.Range("D18").Value = "59.838.70"
$ws.Range("E18").Value = "  -5.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.47"
$ws.Range("E19").Value = "  -1.95%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "423.97"
$ws.Range("E20").Value = "  -8.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.05"
$ws.Range("E21").Value = "  -6.75%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.664"
$ws.Range("E22").Value = "  -4.87%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.93"
$ws.Range("E23").Value = "  -9.83%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.73"
$ws.Range("E24").Value = "  -4.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "78.79"
$ws.Range("E25").Value = "  -5.34%  "

$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("E27").Value = "  +0.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.50"
$ws.Range("E28").Value = "  -7.66%  "

$ws.Range("E29").Value = "  -8.31%  "

$ws.Range("E30").Value = "  -7.80%  "

$ws.Range("E31").Value = "  -11.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "24.92"
$ws.Range("E32").Value = "  -8.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0953"
$ws.Range("E33").Value = "  -5.58%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.52"
$ws.Range("E34").Value = "  -6.42%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.11"
$ws.Range("E35").Value = "  -2.19%  "

$ws.Range("E36").Value = "  -10.19%  "

$ws.Range("E37").Value = "  -19.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.44"
$ws.Range("E38").Value = "  +3.69%  "

$ws.Range("D39").Value = "0.0₃0638"
$ws.Range("E39").Value = "  -13.07%  "

$ws.Range("E40").Value = "  -10.06%  "

$ws.Range("E41").Value = "  -6.18%  "

$ws.Range("D42").Value = "2.651.50"
$ws.Range("E42").Value = "  -5.60%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "362.30"
$ws.Range("E43").Value = "  -8.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.38"
$ws.Range("E44").Value = "  -9.71%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "120.12"
$ws.Range("E46").Value = "  -4.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.230"
$ws.Range("E47").Value = "  -8.54%  "

$ws.Range("E48").Value = "  -4.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.94"
$ws.Range("E49").Value = "  -8.41%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.01"
$ws.Range("E50").Value = "  -8.70%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.98"
$ws.Range("E51").Value = "  -9.22%  "
